$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Standardize the Increase/Decrease labels in column I
$ws.Range("I6").Value = "Increase"
$ws.Range("I7").Value = "Increase"
$ws.Range("I8").Value = "Increase"
$ws.Range("I9").Value = "Decrease"
$ws.Range("I10").Value = "Decrease"

# New (blank) row of test case data, matching the style of the row above
$ws.Range("I10").Copy()
$ws.Range("I11").PasteSpecial(-4122)
$ws.Range("I11").Value = $null

$ws.Range("H15").Select()
